$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.496614933013916
$ws.Range("C1").Value = 1.899153590202332
$ws.Range("D1").Value = 1.488698244094849
$ws.Range("E1").Value = 1.353420376777649
